$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.9022285165660638
$ws.Range("D2").Value = 1567380.082237158
$ws.Range("E2").Value = 1251.950511097447
$ws.Range("F2").Value = 808.5960313687351
$ws.Range("G2").Value = 0.419342489824613

$ws.Range("C3").Value = 0.9509447274214927
$ws.Range("D3").Value = 786407.7997775306
$ws.Range("E3").Value = 886.7963688342046
$ws.Range("F3").Value = 421.6457174638487
$ws.Range("G3").Value = 0.1150866456657527

$ws.Range("C4").Value = 0.9602309778394038
$ws.Range("D4").Value = 637539.4034671161
$ws.Range("E4").Value = 798.4606461605456
$ws.Range("F4").Value = 428.5298000221808
$ws.Range("G4").Value = 0.1162490636065167

$ws.Range("C5").Value = 0.9796780128047003
$ws.Range("D5").Value = 325782.9056354016
$ws.Range("E5").Value = 570.7739531858489
$ws.Range("F5").Value = 278.5198182100496
$ws.Range("G5").Value = 0.06649046156857015

$ws.Range("C6").Value = 0.9807914039581893
$ws.Range("D6").Value = 307934.0702037818
$ws.Range("E6").Value = 554.9180752181188
$ws.Range("F6").Value = 284.1876959300369
$ws.Range("G6").Value = 0.07538644274443188

$ws.Range("C7").Value = 0.9744326434005663
$ws.Range("D7").Value = 409871.7139388064
$ws.Range("E7").Value = 640.212241322209
$ws.Range("F7").Value = 325.8250850328128
$ws.Range("G7").Value = 0.09006429925956252

$ws.Range("C8").Value = 0.9753317332311658
$ws.Range("D8").Value = 395458.3549190896
$ws.Range("E8").Value = 628.8547963712208
$ws.Range("F8").Value = 360.0893533603586
$ws.Range("G8").Value = 0.1210643723634061
